$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4290.361796875724
$ws.Range("C3").Value = 4264.722140426434
$ws.Range("C4").Value = 4182.492051348102
$ws.Range("C5").Value = 4182.492051348102
$ws.Range("C6").Value = 4134.834936170226
$ws.Range("C7").Value = 4112.378057382037
$ws.Range("C8").Value = 4112.378057382037
$ws.Range("C9").Value = 4110.918449915516
$ws.Range("C10").Value = 4059.340755360652
$ws.Range("C11").Value = 3888.828635297548
$ws.Range("C12").Value = 3888.828635297548
